$p = $ppt.ActivePresentation
$s2 = $p.Slides.Item(2)

# --- The new shape needs to land on XML id="100". The engine hands out new
# shape ids from a simple per-slide counter that starts at 2 and silently
# skips over any id already used elsewhere in the deck (ids 2-91 are used by
# the slide master/layouts/slide 1, ids 92-99 are used by slide 2's existing
# shapes). So we burn through the 2..91 range with disposable textboxes,
# create our real shape next (which lands on the first free id, 100), then
# remove the disposable ones again.
$burned = New-Object System.Collections.ArrayList
$n = 1
while ($n -le 90) {
    [void]$burned.Add($s2.Shapes.AddTextbox(1, 0, 0, 10, 10))
    $n = $n + 1
}

# Duplicate the slide's existing body-copy textbox so the new shape inherits
# the exact same spPr/bodyPr/pPr/rPr scaffolding (noFill, no line, Calibri
# run fonts, dk1 solid fill, etc.) that Google Slides originally emitted.
$template = $s2.Shapes.Item(1)
$range = $template.Duplicate()
$shp = $range.Item(1)

foreach ($b in $burned) {
    $b.Delete()
}

$shp.Name = "Google Shape;100;p2"

# Shape resizes to fit its text, matching <a:spAutoFit/>.
$shp.TextFrame.AutoSize = 1

$tr = $shp.TextFrame.TextRange
$tr.Text = "Sentiment measures the intensity of the passenger" + [char]0x2019 + "s review, ranging from -1 (strongly negative) to 1 (strongly positive)."
$tr.Font.Size = 13
$tr.ParagraphFormat.Alignment = 2

$tr.Characters(1, 9).Font.Bold = $true

# Pin the position/size to the exact EMU values from the target deck.
$shp.Left = 397.1496132992126
$shp.Top = 135.9085009370079
$shp.Width = 452.10237220472436
$shp.Height = 46.062994125984254
